# Add 2022-Q4 data
# 1) Insert a new worksheet named "2022-Q4" right before the existing "2022-Q3"
#    worksheet, cloning its layout/format, then fill in the new quarter's
#    fund-holding figures.
# 2) Update the "总计" (totals) summary sheet: shift the quarter rows down by
#    one and insert the new 2022-Q4 figures at the top, adding a trailing row
#    for the 2021-Q2 entry that falls off the bottom.

$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q4" worksheet -----------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# D:G hold numeric-looking figures but are stored as plain text in the
# source file (no "s" style override). Temporarily force a text number
# format so the assignments below aren't auto-coerced to numbers, then
# clear the format again so the cells end up with no style override,
# matching the original sheet's formatting exactly.
$textCells = $q4.Range("D2:G3")
$textCells.NumberFormat = "@"

# Row 2 -> fund 008763
$q4.Range("D2").Value = "21.47"
$q4.Range("E2").Value = "92.63"
$q4.Range("F2").Value = "6.71"
$q4.Range("G2").Value = "1.4406"
$q4.Range("H2").Value = 5

# Row 3 -> fund 008764
$q4.Range("D3").Value = "16.17"
$q4.Range("E3").Value = "92.63"
$q4.Range("F3").Value = "6.71"
$q4.Range("G3").Value = "1.0850"
$q4.Range("H3").Value = 5

$textCells.ClearFormats()

# --- Step 2: update the "总计" summary sheet --------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room for a new row 5 (2021-Q2), copying the formatting of row 4 so the
# new row's A/style attributes match the rest of the table.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 2.28

$total.Range("B4").Value = "2022-Q2"
$total.Range("D4").Value = 2.31

$total.Range("B3").Value = "2022-Q3"
$total.Range("D3").Value = 2.11

$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 2.53

# --- Step 3: restore the originally-selected tab (last sheet, 2021-Q2) -----
$wb.Worksheets.Item("2021-Q2").Activate()
